# Loan RBI, Variable Instalments
# The "Repayment Schedule" sheet gains a new (blank) column inserted
# before the existing "Late" column, pushing "Late" / heading / "Outstanding"
# one column to the right (N -> O -> P -> Q). The sheet also becomes the
# active/selected sheet in the workbook, with the cursor resting on M13.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment Schedule")

# Make "Repayment Schedule" the active sheet (this also clears
# tabSelected on whichever sheet was active before, e.g. NewLoanInput).
$ws.Activate()

# Insert a new blank column before column N, shifting the "Late",
# heading and "Outstanding" columns one place to the right.
$ws.Columns("N").Insert()

# Leave the selection on M13 as recorded in the saved view state.
$ws.Range("M13").Select()
